$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 2.116619666666667
$ws.Range("N2").Value = 6.349859
$ws.Range("O2").Value = 0.02062449887850904
$ws.Range("P2").Value = 0.02062449887850904
$ws.Range("Q2").Value = 0.1260588119477778
$ws.Range("R2").Value = 1.13452930753
$ws.Range("S2").Value = 0.02062449887850904
$ws.Range("T2").Value = 0.02062449887850904

# Row 3 updates
$ws.Range("O3").Value = 0.7564347981517648
$ws.Range("P3").Value = 0.7564347981517647
$ws.Range("S3").Value = 0.7564347981517648
$ws.Range("T3").Value = 0.7564347981517647

# Row 4 updates
$ws.Range("O4").Value = 0.2229407029697262
$ws.Range("P4").Value = 0.2229407029697262
$ws.Range("S4").Value = 0.2229407029697262
$ws.Range("T4").Value = 0.2229407029697262
